# Updated cryptos list on Sat Jun 22 22:26:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (column D) / Volume(1h) (column E) text.
# $null means "leave this column unchanged for this row".
$updates = @(
    @{ Row = 2;  D = "64.276.04"; E = "  +0.15%  " },
    @{ Row = 3;  D = "3.504.54";  E = "  -0.58%  " },
    @{ Row = 4;  D = $null;       E = "  +0.02%  " },
    @{ Row = 5;  D = "588.64";    E = "  +0.20%  " },
    @{ Row = 6;  D = "134.19";    E = "  +0.21%  " },
    @{ Row = 7;  D = $null;       E = "  -0.01%  " },
    @{ Row = 8;  D = $null;       E = "  -0.46%  " },
    @{ Row = 9;  D = $null;       E = "  +0.05%  " },
    @{ Row = 10; D = "7.30";      E = "  +2.11%  " },
    @{ Row = 11; D = "0.386";     E = "  +2.21%  " },
    @{ Row = 12; D = "4.099.10";  E = "  -0.48%  " },
    @{ Row = 13; D = $null;       E = "  +1.16%  " },
    @{ Row = 14; D = $null;       E = "  +0.90%  " },
    @{ Row = 15; D = "3.502.59";  E = "  -0.14%  " },
    @{ Row = 16; D = "64.301.72"; E = "  +0.17%  " },
    @{ Row = 18; D = $null;       E = "  +0.37%  " },
    @{ Row = 19; D = "5.75";      E = "  +2.41%  " },
    @{ Row = 20; D = $null;       E = "  -2.83%  " },
    @{ Row = 21; D = "392.80";    E = "  +2.62%  " },
    @{ Row = 22; D = "0.571";     E = "  -0.08%  " },
    @{ Row = 23; D = "3.643.26";  E = "  -0.54%  " },
    @{ Row = 24; D = "74.64";     E = $null },
    @{ Row = 25; D = "1.00";      E = "  +0.05%  " },
    @{ Row = 26; D = $null;       E = "  -0.41%  " },
    @{ Row = 27; D = $null;       E = "  +0.11%  " },
    @{ Row = 28; D = "7.34";      E = "  -1.86%  " },
    @{ Row = 29; D = $null;       E = "  +0.81%  " },
    @{ Row = 30; D = $null;       E = "  -2.65%  " },
    @{ Row = 31; D = $null;       E = "  -7.62%  " },
    @{ Row = 32; D = "3.525.93";  E = "  -0.26%  " },
    @{ Row = 33; D = "0.154";     E = "  +5.76%  " },
    @{ Row = 34; D = $null;       E = "  +0.05%  " },
    @{ Row = 35; D = "23.43";     E = "  -0.85%  " },
    @{ Row = 36; D = $null;       E = "  -5.01%  " },
    @{ Row = 37; D = "6.86";      E = "  -1.16%  " },
    @{ Row = 38; D = "167.54";    E = "  +4.62%  " },
    @{ Row = 39; D = $null;       E = "  -1.11%  " },
    @{ Row = 40; D = "0.0779";    E = "  -0.93%  " },
    @{ Row = 41; D = "0.810";     E = "  -0.40%  " },
    @{ Row = 42; D = $null;       E = "  +0.05%  " },
    @{ Row = 43; D = "25.26";     E = "  -5.12%  " },
    @{ Row = 44; D = $null;       E = "  -0.62%  " },
    @{ Row = 45; D = $null;       E = "  +2.56%  " },
    @{ Row = 46; D = $null;       E = "  -4.49%  " },
    @{ Row = 47; D = $null;       E = "  -0.87%  " },
    @{ Row = 48; D = "0.891";     E = "  -2.03%  " },
    @{ Row = 49; D = "2.314.81";  E = "  -6.59%  " },
    @{ Row = 50; D = "0.0258";    E = "  -1.72%  " },
    @{ Row = 51; D = "21.15";     E = "  -1.37%  " }
)

# Column D ("Price") holds text such as "64.277.92" or "134.32". Some of the
# new values (e.g. "588.64") look like ordinary decimals, so Excel would
# otherwise silently convert them to numbers on assignment. Force the whole
# Price column to Text first, write the values, then restore the original
# (default/no explicit style) look so the workbook's styling is untouched.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

$priceRange.Style = "Normal"
